$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header columns (descriptive Spanish headers -> short machine-friendly names) ---
$ws.Cells.Item(1,1).Value = "mx_state"
$ws.Cells.Item(1,2).Value = "mx_municipality"
$ws.Cells.Item(1,3).Value = "n_matriculas"
$ws.Cells.Item(1,4).Value = "pct_matriculas"

# --- 2) Fix two floating point percentage values (last-digit precision refresh) ---
$ws.Cells.Item(332,4).Value = 0.009368111682586331
$ws.Cells.Item(500,4).Value = 0.0911094783247612

# --- 3) Title-case the Spanish lowercase connector words (de/del/la/las/los/el/y)
#        in state/municipality names (e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga"),
#        plus one stray double-capitalization fix (MonteMorelos -> Montemorelos).
$renames = @(
    @{R=5; C=2; V="Pabellón De Arteaga"},
    @{R=6; C=2; V="Rincón De Romos"},
    @{R=10; C=2; V="Playas De Rosarito"},
    @{R=23; C=2; V="Amatenango De La Frontera"},
    @{R=28; C=2; V="Chiapa De Corzo"},
    @{R=32; C=2; V="Comitán De Domínguez"},
    @{R=47; C=2; V="Marqués De Comillas"},
    @{R=48; C=2; V="Mazapa De Madero"},
    @{R=51; C=2; V="Ocozocoautla De Espinosa"},
    @{R=57; C=2; V="San Cristóbal De Las Casas"},
    @{R=82; C=2; V="Guadalupe Y Calvo"},
    @{R=83; C=2; V="Hidalgo Del Parral"},
    @{R=90; C=2; V="San Francisco Del Oro"},
    @{R=112; C=1; V="Ciudad De México"},
    @{R=116; C=2; V="Cuajimalpa De Morelos"},
    @{R=131; C=2; V="Coneto De Comonfort"},
    @{R=150; C=2; V="San Juan Del Río"},
    @{R=159; C=1; V="Estado De México"},
    @{R=159; C=2; V="Acambay De Ruíz Castañeda"},
    @{R=160; C=2; V="Almoloya De Juárez"},
    @{R=165; C=2; V="Atizapán De Zaragoza"},
    @{R=170; C=2; V="Coacalco De Berriozábal"},
    @{R=174; C=2; V="Ecatepec De Morelos"},
    @{R=176; C=2; V="Ixtapan De La Sal"},
    @{R=182; C=2; V="Naucalpan De Juárez"},
    @{R=185; C=2; V="San Felipe Del Progreso"},
    @{R=193; C=2; V="Tenango Del Valle"},
    @{R=200; C=2; V="Tlalnepantla De Baz"},
    @{R=204; C=2; V="Valle De Chalco Solidaridad"},
    @{R=205; C=2; V="Villa De Allende"},
    @{R=214; C=2; V="Apaseo El Alto"},
    @{R=215; C=2; V="Apaseo El Grande"},
    @{R=221; C=2; V="Dolores Hidalgo Cuna De La Independencia Nacional"},
    @{R=224; C=2; V="Jaral Del Progreso"},
    @{R=233; C=2; V="San Diego De La Unión"},
    @{R=235; C=2; V="San Francisco Del Rincón"},
    @{R=237; C=2; V="San Luis De La Paz"},
    @{R=238; C=2; V="Santa Cruz De Juventino Rosas"},
    @{R=240; C=2; V="Silao De La Victoria"},
    @{R=245; C=2; V="Valle De Santiago"},
    @{R=250; C=2; V="Acapulco De Juárez"},
    @{R=253; C=2; V="Ajuchitlán Del Progreso"},
    @{R=254; C=2; V="Alcozauca De Guerrero"},
    @{R=258; C=2; V="Atenango Del Río"},
    @{R=260; C=2; V="Atoyac De Álvarez"},
    @{R=261; C=2; V="Ayutla De Los Libres"},
    @{R=263; C=2; V="Buenavista De Cuéllar"},
    @{R=264; C=2; V="Chilapa De Álvarez"},
    @{R=265; C=2; V="Chilpancingo De Los Bravo"},
    @{R=270; C=2; V="Coyuca De Benítez"},
    @{R=271; C=2; V="Coyuca De Catalán"},
    @{R=273; C=2; V="Cutzamala De Pinzón"},
    @{R=279; C=2; V="Huitzuco De Los Figueroa"},
    @{R=280; C=2; V="Iguala De La Independencia"},
    @{R=281; C=2; V="Ixcateopan De Cuauhtémoc"},
    @{R=282; C=2; V="Zihuatanejo De Azueta"},
    @{R=286; C=2; V="Mártir De Cuilapan"},
    @{R=297; C=2; V="Taxco De Alarcón"},
    @{R=299; C=2; V="Técpan De Galeana"},
    @{R=301; C=2; V="Tepecoacuilco De Trujano"},
    @{R=302; C=2; V="Tixtla De Guerrero"},
    @{R=305; C=2; V="Tlapa De Comonfort"},
    @{R=314; C=2; V="Agua Blanca De Iturbide"},
    @{R=317; C=2; V="Atotonilco El Grande"},
    @{R=321; C=2; V="Cuautepec De Hinojosa"},
    @{R=324; C=2; V="Huasca De Ocampo"},
    @{R=325; C=2; V="Huejutla De Reyes"},
    @{R=328; C=2; V="Jacala De Ledezma"},
    @{R=333; C=2; V="Mineral Del Monte"},
    @{R=334; C=2; V="Mixquiahuala De Juárez"},
    @{R=335; C=2; V="Molango De Escamilla"},
    @{R=336; C=2; V="Pachuca De Soto"},
    @{R=339; C=2; V="Progreso De Obregón"},
    @{R=343; C=2; V="Santiago De Anaya"},
    @{R=346; C=2; V="Tenango De Doria"},
    @{R=348; C=2; V="Tezontepec De Aldama"},
    @{R=351; C=2; V="Tula De Allende"},
    @{R=352; C=2; V="Tulancingo De Bravo"},
    @{R=353; C=2; V="Zacualtipán De Ángeles"},
    @{R=357; C=2; V="Ahualulco De Mercado"},
    @{R=360; C=2; V="Atotonilco El Alto"},
    @{R=362; C=2; V="Autlán De Navarro"},
    @{R=371; C=2; V="Encarnación De Díaz"},
    @{R=374; C=2; V="Huejuquilla El Alto"},
    @{R=375; C=2; V="Ixtlahuacán Del Río"},
    @{R=380; C=2; V="Lagos De Moreno"},
    @{R=384; C=2; V="Ojuelos De Jalisco"},
    @{R=389; C=2; V="San Cristóbal De La Barranca"},
    @{R=391; C=2; V="San Juan De Los Lagos"},
    @{R=392; C=2; V="San Juanito De Escobedo"},
    @{R=396; C=2; V="San Miguel El Alto"},
    @{R=397; C=2; V="San Sebastián Del Oeste"},
    @{R=398; C=2; V="Santa María De Los Ángeles"},
    @{R=399; C=2; V="Talpa De Allende"},
    @{R=400; C=2; V="Tamazula De Gordiano"},
    @{R=404; C=2; V="Tepatitlán De Morelos"},
    @{R=405; C=2; V="Tlajomulco De Zúñiga"},
    @{R=410; C=2; V="Unión De San Antonio"},
    @{R=411; C=2; V="Yahualica De González Gallo"},
    @{R=412; C=2; V="Zacoalco De Torres"},
    @{R=415; C=2; V="Zapotlán Del Rey"},
    @{R=416; C=2; V="Zapotlán El Grande"},
    @{R=481; C=2; V="Tiquicheo De Nicolás Romero"},
    @{R=510; C=2; V="Jonacatepec De Leandro Valle"},
    @{R=513; C=2; V="Puente De Ixtla"},
    @{R=518; C=2; V="Tetela Del Volcán"},
    @{R=519; C=2; V="Tlaltizapán De Zapata"},
    @{R=524; C=2; V="Zacualpan De Amilpas"},
    @{R=526; C=2; V="Bahía De Banderas"},
    @{R=528; C=2; V="Ixtlán Del Río"},
    @{R=543; C=2; V="Lampazos De Naranjo"},
    @{R=545; C=2; V="Montemorelos"},
    @{R=548; C=2; V="San Nicolás De Los Garza"},
    @{R=554; C=2; V="Acatlán De Pérez Figueroa"},
    @{R=555; C=2; V="Chalcatongo De Hidalgo"},
    @{R=557; C=2; V="Coicoyán De Las Flores"},
    @{R=559; C=2; V="Eloxochitlán De Flores Magón"},
    @{R=560; C=2; V="Heroica Ciudad De Ejutla De Crespo"},
    @{R=561; C=2; V="Heroica Ciudad De Huajuapan De León"},
    @{R=562; C=2; V="Heroica Ciudad De Tlaxiaco"},
    @{R=563; C=2; V="Ixtlán De Juárez"},
    @{R=564; C=2; V="Heroica Ciudad De Juchitán De Zaragoza"},
    @{R=567; C=2; V="Mariscala De Juárez"},
    @{R=568; C=2; V="Miahuatlán De Porfirio Díaz"},
    @{R=570; C=2; V="Nejapa De Madero"},
    @{R=571; C=2; V="Oaxaca De Juárez"},
    @{R=572; C=2; V="Ocotlán De Morelos"},
    @{R=573; C=2; V="Pinotepa De Don Luis"},
    @{R=575; C=2; V="Putla Villa De Guerrero"},
    @{R=576; C=2; V="Reforma De Pineda"},
    @{R=585; C=2; V="San Francisco Del Mar"},
    @{R=647; C=2; V="Santo Domingo De Morelos"},
    @{R=655; C=2; V="Tlacolula De Matamoros"},
    @{R=656; C=2; V="Totontepec Villa De Morelos"},
    @{R=657; C=2; V="Villa De Tututepec De Melchor Ocampo"},
    @{R=658; C=2; V="Villa Sola De Vega"},
    @{R=659; C=2; V="Zimatlán De Álvarez"},
    @{R=667; C=2; V="Ayotoxco De Guerrero"},
    @{R=674; C=2; V="Cuapiaxtla De Madero"},
    @{R=675; C=2; V="Cuayuca De Andrade"},
    @{R=676; C=2; V="Cuetzalan Del Progreso"},
    @{R=687; C=2; V="Ixcamilpa De Guerrero"},
    @{R=689; C=2; V="Izúcar De Matamoros"},
    @{R=692; C=2; V="Los Reyes De Juárez"},
    @{R=703; C=2; V="San Salvador El Verde"},
    @{R=710; C=2; V="Tepatlaxco De Hidalgo"},
    @{R=713; C=2; V="Tepexi De Rodríguez"},
    @{R=724; C=2; V="Xayacatlán De Bravo"},
    @{R=734; C=2; V="Amealco De Bonfil"},
    @{R=735; C=2; V="Cadereyta De Montes"},
    @{R=738; C=2; V="Jalpan De Serra"},
    @{R=739; C=2; V="Landa De Matamoros"},
    @{R=740; C=2; V="Pinal De Amoles"},
    @{R=742; C=2; V="San Juan Del Río"},
    @{R=754; C=2; V="Cerro De San Pedro"},
    @{R=756; C=2; V="Ciudad Del Maíz"},
    @{R=764; C=2; V="Mexquitic De Carmona"},
    @{R=772; C=2; V="Santa María Del Río"},
    @{R=773; C=2; V="Soledad De Graciano Sánchez"},
    @{R=779; C=2; V="Tanquián De Escobedo"},
    @{R=783; C=2; V="Villa De Arriaga"},
    @{R=784; C=2; V="Villa De Guadalupe"},
    @{R=785; C=2; V="Villa De Ramos"},
    @{R=786; C=2; V="Villa De Reyes"},
    @{R=854; C=2; V="Alto Lucero De Gutiérrez Barrios"},
    @{R=857; C=2; V="Amatlán De Los Reyes"},
    @{R=862; C=2; V="Camarón De Tejeda"},
    @{R=863; C=2; V="Castillo De Teayo"},
    @{R=873; C=2; V="Cosamaloapan De Carpio"},
    @{R=882; C=2; V="Hueyapan De Ocampo"},
    @{R=883; C=2; V="Ignacio De La Llave"},
    @{R=886; C=2; V="Ixhuatlán Del Sureste"},
    @{R=892; C=2; V="Juchique De Ferrer"},
    @{R=895; C=2; V="Lerdo De Tejada"},
    @{R=897; C=2; V="Martínez De La Torre"},
    @{R=908; C=2; V="Paso De Ovejas"},
    @{R=910; C=2; V="Poza Rica De Hidalgo"},
    @{R=916; C=2; V="Sayula De Alemán"},
    @{R=935; C=2; V="Vega De Alatorre"},
    @{R=942; C=2; V="Zontecomatlán De López Y Fuentes"},
    @{R=952; C=2; V="Concepción Del Oro"},
    @{R=953; C=2; V="El Plateado De Joaquín Amaro"},
    @{R=966; C=2; V="Moyahua De Estrada"},
    @{R=967; C=2; V="Noria De Ángeles"},
    @{R=973; C=2; V="Tlaltenango De Sánchez Román"},
    @{R=975; C=2; V="Villa De Cos"}
)

foreach ($item in $renames) {
    $ws.Cells.Item($item.R, $item.C).Value = $item.V
}

# --- 4) Drop the trailing metadata/footer rows (984-988): sample size, source, author, date ---
$ws.Range("A984:A988").EntireRow.Delete() | Out-Null
